# cap nhat file Anh12.xlsx
# Fill in the "D" (correct/incorrect) boolean column on the "DapAn" sheet for
# every question row that was still missing it, marking them as FALSE, and
# leave the view selecting the first block of newly-edited cells (D3:D7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DapAn")
$ws.Activate()

# Contiguous row blocks (on column D) that had no value yet and need FALSE.
$ranges = @(
    "D3:D7",
    "D9:D10",
    "D12:D13",
    "D15:D18",
    "D20:D22",
    "D24:D28",
    "D30:D32",
    "D34:D35",
    "D37",
    "D39:D42",
    "D44:D47",
    "D49",
    "D51:D53",
    "D55:D57",
    "D59:D62",
    "D64:D66",
    "D68:D71",
    "D73:D76",
    "D79:D81"
)

foreach ($rng in $ranges) {
    $ws.Range($rng).Value = $false
}

# Match the saved selection/view state: active cell D3, selected block D3:D7,
# scrolled back to the top of the sheet.
$ws.Range("D3:D7").Select()
